$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage (preserves exact
# string formatting such as trailing zeros / multi-dot price strings) and
# then restores the cells original style so no formatting diff is left
# behind.
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '67.347.57'
$ws.Range('E2').Value = '  +0.62%  '
Set-TextValue $ws.Range('D3') '3.523.35'
$ws.Range('E3').Value = '  +0.65%  '
$ws.Range('E4').Value = '  +0.00%  '
Set-TextValue $ws.Range('D5') '596.88'
$ws.Range('E5').Value = '  +0.63%  '
Set-TextValue $ws.Range('D6') '173.75'
$ws.Range('E6').Value = '  +2.83%  '
Set-TextValue $ws.Range('D7') '1.00'
$ws.Range('E8').Value = '  +3.35%  '
$ws.Range('E9').Value = '  +8.23%  '
$ws.Range('E10').Value = '  +0.31%  '
Set-TextValue $ws.Range('D11') '0.438'
$ws.Range('E11').Value = '  -0.33%  '
Set-TextValue $ws.Range('D12') '4.133.78'
$ws.Range('E12').Value = '  +0.62%  '
$ws.Range('E13').Value = '  +0.01%  '
Set-TextValue $ws.Range('D14') '28.81'
$ws.Range('E14').Value = '  +2.25%  '
Set-TextValue $ws.Range('D15') '0.0000183'
$ws.Range('E15').Value = '  +1.87%  '
Set-TextValue $ws.Range('D16') '67.284.72'
$ws.Range('E16').Value = '  +0.62%  '
Set-TextValue $ws.Range('D17') '3.496.62'
$ws.Range('E17').Value = '  -0.33%  '
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('E19').Value = '  +1.59%  '
Set-TextValue $ws.Range('D20') '398.06'
$ws.Range('E20').Value = '  +0.91%  '
Set-TextValue $ws.Range('D21') '8.00'
$ws.Range('E21').Value = '  +0.83%  '
Set-TextValue $ws.Range('D22') '73.49'
$ws.Range('E22').Value = '  +0.01%  '
Set-TextValue $ws.Range('D23') '0.540'
$ws.Range('E23').Value = '  +2.03%  '
Set-TextValue $ws.Range('D24') '0.998'
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('E25').Value = '  -3.35%  '
Set-TextValue $ws.Range('D26') '10.29'
$ws.Range('E26').Value = '  +2.70%  '
Set-TextValue $ws.Range('D27') '0.181'
$ws.Range('E27').Value = '  -0.69%  '
$ws.Range('E28').Value = '  -0.16%  '
Set-TextValue $ws.Range('D29') '6.28'
$ws.Range('E29').Value = '  -1.65%  '
$ws.Range('E30').Value = '  -0.47%  '
Set-TextValue $ws.Range('D31') '2.08'
$ws.Range('E31').Value = '  +1.02%  '
Set-TextValue $ws.Range('D32') '24.17'
$ws.Range('E32').Value = '  +2.65%  '
Set-TextValue $ws.Range('D33') '7.41'
$ws.Range('E33').Value = '  -0.30%  '
$ws.Range('E34').Value = '  +2.30%  '
Set-TextValue $ws.Range('D35') '163.63'
$ws.Range('E35').Value = '  +0.89%  '
Set-TextValue $ws.Range('D36') '0.897'
$ws.Range('E36').Value = '  -0.68%  '
$ws.Range('E37').Value = '  -1.06%  '
Set-TextValue $ws.Range('D38') '6.94'
$ws.Range('E38').Value = '  +3.63%  '
$ws.Range('E39').Value = '  +0.97%  '
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D40') '27.55'
$ws.Range('E40').Value = '  +3.77%  '
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D41') '0.0746'
$ws.Range('E41').Value = '  -0.84%  '
$ws.Range('E42').Value = '  +0.37%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D43') '2.63'
$ws.Range('E43').Value = '  +3.62%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D44') '2.811.22'
$ws.Range('E44').Value = '  -0.80%  '
Set-TextValue $ws.Range('D45') '42.93'
$ws.Range('E45').Value = '  -1.32%  '
$ws.Range('E46').Value = '  -2.33%  '
Set-TextValue $ws.Range('D47') '340.80'
$ws.Range('E47').Value = '  -2.31%  '
$ws.Range('E48').Value = '  +1.83%  '
Set-TextValue $ws.Range('D49') '33.65'
$ws.Range('E49').Value = '  +0.36%  '
$ws.Range('E50').Value = '  +0.34%  '
Set-TextValue $ws.Range('D51') '0.852'
$ws.Range('E51').Value = '  -0.33%  '
